$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DSD")

# 1) REF_AREA -> GEO_PICT (row 4, column A). Rest of the row (B..G) stays the same.
$ws.Range("A4").Value = "GEO_PICT"

# 2) UNIT_MEASURE row (row 14): codelist reference CL_UNIT_MEASURE -> CL_COM_UNIT_MEASURE
$ws.Range("F14").Value = "CL_COM_UNIT_MEASURE"

# 3) Insert a new row at 15 for UNIT_MULT (pushes OBS_STATUS/COMMENT down by one)
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "UNIT_MULT"
$ws.Range("B15").Value = "Unit multiplier"
$ws.Range("C15").Value = "Y"
$ws.Range("D15").Value = "Attribute"
$ws.Range("E15").Value = "Coded"
$ws.Range("F15").Value = "CL_COM_UNIT_MULT"
$ws.Range("G15").Value = "Y"

# OBS_STATUS is now on row 16 (shifted down by the insert above): update its codelist
$ws.Range("F16").Value = "CL_COM_OBS_STATUS"

# 4) Insert a new row at 17 for DATA_SOURCE (pushes the old COMMENT row down to 18)
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "DATA_SOURCE"
$ws.Range("B17").Value = "Data source"
$ws.Range("C17").Value = "N"
$ws.Range("D17").Value = "Attribute"
$ws.Range("E17").Value = "Uncoded"
$ws.Range("F17").Value = "Text"
$ws.Range("G17").Value = "N"

# COMMENT row is now on row 18: rename concept id to OBS_COMMENT (label/other cols unchanged)
$ws.Range("A18").Value = "OBS_COMMENT"

# 5) Append new row 19 for CONF_STATUS
$ws.Range("A19").Value = "CONF_STATUS"
$ws.Range("B19").Value = "Confidentiality status"
$ws.Range("C19").Value = "Y"
$ws.Range("D19").Value = "Attribute"
$ws.Range("E19").Value = "Coded"
$ws.Range("F19").Value = "CL_COM_CONF_STATUS"
$ws.Range("G19").Value = "Y"

# Update the sheet's selection to match the new layout (A14, extended to the new data block)
$ws.Range("A14:G19").Select()
